$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:B5 with new values
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 106

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 75

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 62

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 55

# Remove row 6 entirely (was A6=1, B6=52)
$ws.Range("A6:B6").Delete()
